$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.104.97"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.790.07"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0691"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "2.047.82"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.46%  "
$ws.Range("D14").Value = "1.788.43"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "34.094.88"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "1.421.68"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.25%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.921"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.60%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0508"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "1.949.10"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  +0.12%  "
